$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-BOM")

# Fill in the previously blank "Value", "Type" and "System" columns
# for the two connector rows (row 2: MCX connector, row 3: RP-SMA connector).
# A leading apostrophe forces Excel to treat the entry as explicit text
# (quote-prefixed), matching the formatting already used by the other
# text cells in this sheet.

$ws.Range("F2").Value = "'MCX-F"
$ws.Range("H2").Value = "'SMD"
$ws.Range("K2").Value = "'Connector"

$ws.Range("F3").Value = "'RP-SMA"
$ws.Range("H3").Value = "'SMD"
$ws.Range("K3").Value = "'Connector"
